$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.198.54'
$ws.Range("E2").Value = '  -1.68%  '
$ws.Range("D3").Value = '1.876.14'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("A100").Value = "'1.001"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("A100").Value = "'235.72"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E5").Value = '  -1.47%  '
$ws.Range("A100").Value = "'1.000"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("A100").Value = "'0.4835"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E7").Value = '  -1.03%  '
$ws.Range("A100").Value = "'0.2865"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E8").Value = '  -3.65%  '
$ws.Range("A100").Value = "'0.06580"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E9").Value = '  -2.33%  '
$ws.Range("D10").Value = '1.867.70'
$ws.Range("E10").Value = '  -1.08%  '
$ws.Range("A100").Value = "'16.67"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E11").Value = '  -1.46%  '
$ws.Range("A100").Value = "'0.07275"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("A100").Value = "'5.142"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("A100").Value = "'86.89"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E14").Value = '  -3.09%  '
$ws.Range("A100").Value = "'0.6524"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E15").Value = '  -2.15%  '
$ws.Range("D16").Value = '30.180.64'
$ws.Range("E16").Value = '  -1.61%  '
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("A100").Value = "'13.28"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("A100").Value = "'0.9997"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("A100").Value = "'0.000007689"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E19").Value = '  -3.25%  '
$ws.Range("D20").Value = '2.113.60'
$ws.Range("E20").Value = '  -0.60%  '
$ws.Range("A100").Value = "'5.314"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E21").Value = '  +7.39%  '
$ws.Range("A100").Value = "'1.001"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("A100").Value = "'195.49"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E23").Value = '  -6.28%  '
$ws.Range("A100").Value = "'6.098"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("A100").Value = "'9.319"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E25").Value = '  -3.32%  '
$ws.Range("A100").Value = "'160.36"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E26").Value = '  +1.47%  '
$ws.Range("A100").Value = "'18.08"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E27").Value = '  -4.50%  '
$ws.Range("A100").Value = "'1.911"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E28").Value = '  +2.08%  '
$ws.Range("A100").Value = "'1.440"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E29").Value = '  +1.49%  '
$ws.Range("A100").Value = "'4.259"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E30").Value = '  -1.76%  '
$ws.Range("A100").Value = "'0.09099"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E31").Value = '  -0.37%  '
$ws.Range("A100").Value = "'4.052"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("A100").Value = "'0.05126"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E33").Value = '  -0.97%  '
$ws.Range("A100").Value = "'0.7199"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E34").Value = '  -4.78%  '
$ws.Range("A100").Value = "'1.094"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E35").Value = '  -1.10%  '
$ws.Range("A100").Value = "'2.694"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("A100").Value = "'0.01793"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E37").Value = '  -2.33%  '
$ws.Range("A100").Value = "'2.634"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E38").Value = '  -2.26%  '
$ws.Range("A100").Value = "'0.9148"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("A100").Value = "'2.032"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E40").Value = '  -3.12%  '
$ws.Range("A100").Value = "'105.72"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E41").Value = '  -0.83%  '
$ws.Range("A100").Value = "'0.4268"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E42").Value = '  -5.48%  '
$ws.Range("A100").Value = "'5.787"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("A100").Value = "'0.9984"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("A100").Value = "'66.64"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E45").Value = '  +2.53%  '
$ws.Range("A100").Value = "'7.385"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E46").Value = '  -6.36%  '
$ws.Range("E47").Value = '  -3.22%  '
$ws.Range("A100").Value = "'9.052"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E48").Value = '  +1.59%  '
$ws.Range("A100").Value = "'0.05754"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E49").Value = '  -2.38%  '
$ws.Range("A100").Value = "'33.79"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E50").Value = '  -3.74%  '
$ws.Range("A100").Value = "'0.3819"
$ws.Range("A100").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163, $null) | Out-Null
$ws.Range("E51").Value = '  -6.80%  '
$ws.Range("A100").Clear() | Out-Null
$excel.CutCopyMode = $false
